# Applies the "Added updated xlsx to find model that gives min metric" edit.
#
# For each metric column pair on the Metrics sheet (No-Preprocessing in
# C:F, With-Preprocessing in I:L), write a dynamic-array CELL/INDEX/MATCH
# formula into row 12 that reports the address of the row (within 5:10)
# holding the MIN value of that column - i.e. which model gives the best
# (lowest) value for that metric.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metrics")

# Remove the stale "_xlchart.v1.*" hidden defined names left over from a
# previous chart-data edit; Excel regenerates these as needed and the
# workbook no longer references them.
while ($wb.Names.Count -gt 0) {
    $wb.Names.Item(1).Delete()
}

$cols = @("C", "D", "E", "F", "I", "J", "K", "L")
foreach ($col in $cols) {
    $formula = '=CELL("address",INDEX(' + $col + '5:' + $col + '10,MATCH(MIN(' + $col + '5:' + $col + '10),' + $col + '5:' + $col + '10,0)))'
    $ws.Range($col + "12").Formula2 = $formula
}

$excel.Calculate()

# Make "Metrics" the active sheet/tab with H12 as the selected cell.
$ws.Activate() | Out-Null
$ws.Range("H12").Select() | Out-Null
